# Weekly fruit/vegetable price update: a new sample (week of 2022-08-30)
# is inserted for "Feria Lagunitas de Puerto Montt - Acelga" right before
# the existing row 94, pushing the old rows 94-209 down to 95-210.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a brand-new row at position 94 (shifts old row 94..209 -> 95..210).
$ws.Rows.Item(94).Insert()

# Populate the newly-inserted row 94 with this week's data point.
$ws.Range("A94").Value = 4
$ws.Range("B94").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C94").Value = "Los Lagos"
$ws.Range("D94").Value = 44803
$ws.Range("E94").Value = 10
$ws.Range("F94").Value = 100112009
$ws.Range("G94").Value = "Acelga"
$ws.Range("H94").Value = "Sin especificar"
$ws.Range("I94").Value = "Primera"
$ws.Range("J94").Value = 150
$ws.Range("K94").Value = 3000
$ws.Range("L94").Value = 3000
$ws.Range("M94").Value = 3000
$ws.Range("N94").Value = "`$/docena de atados (4 kilos)"
$ws.Range("O94").Value = "Región del Maule"
$ws.Range("P94").Value = 750
$ws.Range("Q94").Value = 4
$ws.Range("R94").Value = "Hortaliza"
